$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.227772
$ws.Range("H2").Value = 36.683316
$ws.Range("I2").Value = 0.6223179025846677
$ws.Range("J2").Value = 0.6245980896688198
$ws.Range("M2").Value = 38.10639333333333
$ws.Range("N2").Value = 114.31918
$ws.Range("O2").Value = 0.3831479157160237
$ws.Range("P2").Value = 0.4159903984418967
$ws.Range("Q2").Value = 465.95628942232
$ws.Range("R2").Value = 4193.60660480088
$ws.Range("S2").Value = 0.2384398072880829
$ws.Range("T2").Value = 0.2598268081873799
$ws.Range("G3").Value = 12.227772
$ws.Range("H3").Value = 36.683316
$ws.Range("I3").Value = 0.6223179025846677
$ws.Range("J3").Value = 0.6245980896688198
$ws.Range("O3").Value = 0.09199521176963764
$ws.Range("P3").Value = 0.09988081163714851
$ws.Range("Q3").Value = 111.877804269648
$ws.Range("R3").Value = 1006.900238426832
$ws.Range("S3").Value = 0.05725026723631323
$ws.Range("T3").Value = 0.06238536414313419
$ws.Range("G4").Value = 12.227772
$ws.Range("H4").Value = 36.683316
$ws.Range("I4").Value = 0.6223179025846677
$ws.Range("J4").Value = 0.6245980896688198
$ws.Range("M4").Value = 15.023598
$ws.Range("N4").Value = 45.070794
$ws.Range("O4").Value = 0.1510575983904562
$ws.Range("P4").Value = 0.1640058785774412
$ws.Range("Q4").Value = 183.705130963656
$ws.Range("R4").Value = 1653.346178672904
$ws.Range("S4").Value = 0.09400584779982576
$ws.Range("T4").Value = 0.1024377584539262
$ws.Range("G5").Value = 12.227772
$ws.Range("H5").Value = 36.683316
$ws.Range("I5").Value = 0.6223179025846677
$ws.Range("J5").Value = 0.6245980896688198
$ws.Range("M5").Value = 23.556204
$ws.Range("N5").Value = 47.112408
$ws.Range("O5").Value = 0.236850294013169
$ws.Range("P5").Value = 0.1714350065796238
$ws.Range("Q5").Value = 288.039891697488
$ws.Range("R5").Value = 1728.239350184928
$ws.Range("S5").Value = 0.1473961781968372
$ws.Range("T5").Value = 0.1070779776119946
$ws.Range("G6").Value = 12.227772
$ws.Range("H6").Value = 36.683316
$ws.Range("I6").Value = 0.6223179025846677
$ws.Range("J6").Value = 0.6245980896688198
$ws.Range("M6").Value = 13.62041
$ws.Range("N6").Value = 40.86123000000001
$ws.Range("O6").Value = 0.1369489801107134
$ws.Range("P6").Value = 0.1486879047638899
$ws.Range("Q6").Value = 166.54726802652
$ws.Range("R6").Value = 1498.92541223868
$ws.Range("S6").Value = 0.08522580206360852
$ws.Range("T6").Value = 0.09287018127238503
$ws.Range("G7").Value = 3.888411
$ws.Range("H7").Value = 11.665233
$ws.Range("I7").Value = 0.1978960499023984
$ws.Range("J7").Value = 0.1986211455731449
$ws.Range("M7").Value = 38.10639333333333
$ws.Range("N7").Value = 114.31918
$ws.Range("O7").Value = 0.3831479157160237
$ws.Range("P7").Value = 0.4159903984418967
$ws.Range("Q7").Value = 148.17331900766
$ws.Range("R7").Value = 1333.55987106894
$ws.Range("S7").Value = 0.07582345904853817
$ws.Range("T7").Value = 0.08262448948595849
$ws.Range("G8").Value = 3.888411
$ws.Range("H8").Value = 11.665233
$ws.Range("I8").Value = 0.1978960499023984
$ws.Range("J8").Value = 0.1986211455731449
$ws.Range("O8").Value = 0.09199521176963764
$ws.Range("P8").Value = 0.09988081163714851
$ws.Range("Q8").Value = 35.57695422992399
$ws.Range("R8").Value = 320.1925880693159
$ws.Range("S8").Value = 0.01820548901914592
$ws.Range("T8").Value = 0.01983844122814594
$ws.Range("G9").Value = 3.888411
$ws.Range("H9").Value = 11.665233
$ws.Range("I9").Value = 0.1978960499023984
$ws.Range("J9").Value = 0.1986211455731449
$ws.Range("M9").Value = 15.023598
$ws.Range("N9").Value = 45.070794
$ws.Range("O9").Value = 0.1510575983904562
$ws.Range("P9").Value = 0.1640058785774412
$ws.Range("Q9").Value = 58.41792372277799
$ws.Range("R9").Value = 525.7613135050019
$ws.Range("S9").Value = 0.02989370202921417
$ws.Range("T9").Value = 0.03257503548378148
$ws.Range("G10").Value = 3.888411
$ws.Range("H10").Value = 11.665233
$ws.Range("I10").Value = 0.1978960499023984
$ws.Range("J10").Value = 0.1986211455731449
$ws.Range("M10").Value = 23.556204
$ws.Range("N10").Value = 47.112408
$ws.Range("O10").Value = 0.236850294013169
$ws.Range("P10").Value = 0.1714350065796238
$ws.Range("Q10").Value = 91.59620275184399
$ws.Range("R10").Value = 549.577216511064
$ws.Range("S10").Value = 0.04687173760342782
$ws.Range("T10").Value = 0.03405061739818452
$ws.Range("G11").Value = 3.888411
$ws.Range("H11").Value = 11.665233
$ws.Range("I11").Value = 0.1978960499023984
$ws.Range("J11").Value = 0.1986211455731449
$ws.Range("M11").Value = 13.62041
$ws.Range("N11").Value = 40.86123000000001
$ws.Range("O11").Value = 0.1369489801107134
$ws.Range("P11").Value = 0.1486879047638899
$ws.Range("Q11").Value = 52.96175206851
$ws.Range("R11").Value = 476.65576861659
$ws.Range("S11").Value = 0.0271016622020723
$ws.Range("T11").Value = 0.02953256197707447
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.215192
$ws.Range("H12").Value = 0.430384
$ws.Range("I12").Value = 0.01095194072092608
$ws.Range("J12").Value = 0.007328045922130521
$ws.Range("M12").Value = 38.10639333333333
$ws.Range("N12").Value = 114.31918
$ws.Range("O12").Value = 0.3831479157160237
$ws.Range("P12").Value = 0.4159903984418967
$ws.Range("Q12").Value = 8.200190994186666
$ws.Range("R12").Value = 49.20114596511999
$ws.Range("S12").Value = 0.004196213260268276
$ws.Range("T12").Value = 0.003048396742947592
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.215192
$ws.Range("H13").Value = 0.430384
$ws.Range("I13").Value = 0.01095194072092608
$ws.Range("J13").Value = 0.007328045922130521
$ws.Range("O13").Value = 0.09199521176963764
$ws.Range("P13").Value = 0.09988081163714851
$ws.Range("Q13").Value = 1.968895760928
$ws.Range("R13").Value = 11.813374565568
$ws.Range("S13").Value = 0.001007526105910113
$ws.Range("T13").Value = 0.0007319311744166929
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.5
$ws.Range("G14").Value = 0.215192
$ws.Range("H14").Value = 0.430384
$ws.Range("I14").Value = 0.01095194072092608
$ws.Range("J14").Value = 0.007328045922130521
$ws.Range("M14").Value = 15.023598
$ws.Range("N14").Value = 45.070794
$ws.Range("O14").Value = 0.1510575983904562
$ws.Range("P14").Value = 0.1640058785774412
$ws.Range("Q14").Value = 3.232958100816
$ws.Range("R14").Value = 19.397748604896
$ws.Range("S14").Value = 0.001654373863017736
$ws.Range("T14").Value = 0.001201842609714852
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.5
$ws.Range("G15").Value = 0.215192
$ws.Range("H15").Value = 0.430384
$ws.Range("I15").Value = 0.01095194072092608
$ws.Range("J15").Value = 0.007328045922130521
$ws.Range("M15").Value = 23.556204
$ws.Range("N15").Value = 47.112408
$ws.Range("O15").Value = 0.236850294013169
$ws.Range("P15").Value = 0.1714350065796238
$ws.Range("Q15").Value = 5.069106651168
$ws.Range("R15").Value = 20.276426604672
$ws.Range("S15").Value = 0.002593970379766141
$ws.Range("T15").Value = 0.001256283600876232
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.5
$ws.Range("G16").Value = 0.215192
$ws.Range("H16").Value = 0.430384
$ws.Range("I16").Value = 0.01095194072092608
$ws.Range("J16").Value = 0.007328045922130521
$ws.Range("M16").Value = 13.62041
$ws.Range("N16").Value = 40.86123000000001
$ws.Range("O16").Value = 0.1369489801107134
$ws.Range("P16").Value = 0.1486879047638899
$ws.Range("Q16").Value = 2.93100326872
$ws.Range("R16").Value = 17.58601961232
$ws.Range("S16").Value = 0.001499857111963818
$ws.Range("T16").Value = 0.001089591794175154
$ws.Range("G17").Value = 3.31738
$ws.Range("H17").Value = 9.95214
$ws.Range("I17").Value = 0.168834106792008
$ws.Range("J17").Value = 0.1694527188359048
$ws.Range("M17").Value = 38.10639333333333
$ws.Range("N17").Value = 114.31918
$ws.Range("O17").Value = 0.3831479157160237
$ws.Range("P17").Value = 0.4159903984418967
$ws.Range("Q17").Value = 126.4133871161333
$ws.Range("R17").Value = 1137.7204840452
$ws.Range("S17").Value = 0.06468843611913443
$ws.Range("T17").Value = 0.07049070402561072
$ws.Range("G18").Value = 3.31738
$ws.Range("H18").Value = 9.95214
$ws.Range("I18").Value = 0.168834106792008
$ws.Range("J18").Value = 0.1694527188359048
$ws.Range("O18").Value = 0.09199521176963764
$ws.Range("P18").Value = 0.09988081163714851
$ws.Range("Q18").Value = 30.35231523192
$ws.Range("R18").Value = 273.17083708728
$ws.Range("S18").Value = 0.01553192940826839
$ws.Range("T18").Value = 0.01692507509145169
$ws.Range("G19").Value = 3.31738
$ws.Range("H19").Value = 9.95214
$ws.Range("I19").Value = 0.168834106792008
$ws.Range("J19").Value = 0.1694527188359048
$ws.Range("M19").Value = 15.023598
$ws.Range("N19").Value = 45.070794
$ws.Range("O19").Value = 0.1510575983904562
$ws.Range("P19").Value = 0.1640058785774412
$ws.Range("Q19").Value = 49.83898353324
$ws.Range("R19").Value = 448.55085179916
$ws.Range("S19").Value = 0.02550367469839853
$ws.Range("T19").Value = 0.02779124203001869
$ws.Range("G20").Value = 3.31738
$ws.Range("H20").Value = 9.95214
$ws.Range("I20").Value = 0.168834106792008
$ws.Range("J20").Value = 0.1694527188359048
$ws.Range("M20").Value = 23.556204
$ws.Range("N20").Value = 47.112408
$ws.Range("O20").Value = 0.236850294013169
$ws.Range("P20").Value = 0.1714350065796238
$ws.Range("Q20").Value = 78.14488002552001
$ws.Range("R20").Value = 468.86928015312
$ws.Range("S20").Value = 0.03998840783313786
$ws.Range("T20").Value = 0.02905012796856849
$ws.Range("G21").Value = 3.31738
$ws.Range("H21").Value = 9.95214
$ws.Range("I21").Value = 0.168834106792008
$ws.Range("J21").Value = 0.1694527188359048
$ws.Range("M21").Value = 13.62041
$ws.Range("N21").Value = 40.86123000000001
$ws.Range("O21").Value = 0.1369489801107134
$ws.Range("P21").Value = 0.1486879047638899
$ws.Range("Q21").Value = 45.18407572580001
$ws.Range("R21").Value = 406.6566815322001
$ws.Range("S21").Value = 0.02312165873306876
$ws.Range("T21").Value = 0.02519556972025522
